# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> new F value }
$updates = @{
    "展览" = @{
        2  = 5490
        4  = 349
        10 = 327
        11 = 421
        12 = 3005
        14 = 1606
        15 = 14
    }
    "全部类型" = @{
        2  = 5490
        4  = 349
        11 = 327
        12 = 421
        13 = 3005
        15 = 1606
        16 = 14
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
